$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quote row appended at the end of the table (row 41).
# Columns A (date) and B (time) look like date/time literals, so Excel's
# normal smart-parsing would silently coerce them into date/time serials.
# Force them to be stored as literal text (matching the existing rows,
# which are all plain strings) by temporarily marking the cells as Text
# before assigning, then clearing the format again so no stray
# number-format/style ends up attached to the cell.
$ws.Cells.Item(41, 1).NumberFormat = "@"
$ws.Cells.Item(41, 1).Value = "2025-09-26"
$ws.Cells.Item(41, 1).ClearFormats()

$ws.Cells.Item(41, 2).NumberFormat = "@"
$ws.Cells.Item(41, 2).Value = "15:20:15"
$ws.Cells.Item(41, 2).ClearFormats()

# Column C is unambiguously text already ("1.00 EUR = ..."), no coercion risk.
$ws.Cells.Item(41, 3).Value = "1.00 EUR = 1,619.7750"
